$wb = $excel.ActiveWorkbook

# --- 1. "About" sheet: bump the date in C1 (report/version date) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45387

# --- 2. Rename all the "... : NoSettings" labels used as row headers on
#        the "BAU Emissions" sheet to "... : test" ---
$wsBau = $wb.Worksheets.Item("BAU Emissions")
$wsBau.Range("A1:A300").Replace(" : NoSettings", " : test")

# --- 3. Update the data values in row 94 (M:AE) on "BAU Emissions" ---
$wsBau.Range("M94").Value = 1001080
$wsBau.Range("N94").Value = 2002150
$wsBau.Range("O94").Value = 3003230
$wsBau.Range("P94").Value = 4004300
$wsBau.Range("Q94:AE94").Value = 5005380

# --- 4. Sheet view / selection bookkeeping ---
# "Current and Planned Capacity" is no longer the tab shown on open;
# "BAU Emissions" becomes the active / selected sheet+range instead.
$wsCurrent = $wb.Worksheets.Item("Current and Planned Capacity")
$wsCurrent.Activate()

$wsBau.Activate()
$wsBau.Range("A30:AE280").Select()
$excel.ActiveWindow.ScrollRow = 265
$excel.ActiveWindow.ScrollColumn = 1
